# Quiz Results leaderboard grew by one entry: a new completed attempt
# ("Deepa", I-Number I0796921) is inserted as row 6 (ID 112), pushing the
# previously-last two rows ("Meow" ID 109 and "Mona" ID 110) down by one
# row each.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6-7 down to 7-8.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row with the new leaderboard entry.
$ws.Cells.Item(6, 1).Value = 112
$ws.Cells.Item(6, 2).Value = "Deepa"
$ws.Cells.Item(6, 3).Value = "I0796921"
$ws.Cells.Item(6, 4).Value = 60
$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = 3
$ws.Cells.Item(6, 7).Value = 12
$ws.Cells.Item(6, 8).Value = "2025-04-27 23:12:00"
